$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.180165333333334
$ws.Range("H2").Value = 6.540496
$ws.Range("I2").Value = 0.01970539991828544
$ws.Range("J2").Value = 0.01970539991828544
$ws.Range("M2").Value = 28.31444233333334
$ws.Range("N2").Value = 84.94332700000001
$ws.Range("O2").Value = 0.2747173016130739
$ws.Range("P2").Value = 0.2747173016130739
$ws.Range("Q2").Value = 61.73016560779912
$ws.Range("R2").Value = 555.5714904701921
$ws.Range("S2").Value = 0.005413414292757863
$ws.Range("T2").Value = 0.005413414292757863
$ws.Range("G3").Value = 2.180165333333334
$ws.Range("H3").Value = 6.540496
$ws.Range("I3").Value = 0.01970539991828544
$ws.Range("J3").Value = 0.01970539991828544
$ws.Range("O3").Value = 0.2090339131726295
$ws.Range("P3").Value = 0.2090339131726295
$ws.Range("Q3").Value = 46.9708241964569
$ws.Range("R3").Value = 422.737417768112
$ws.Range("S3").Value = 0.004119096855550819
$ws.Range("T3").Value = 0.004119096855550819
$ws.Range("G4").Value = 2.180165333333334
$ws.Range("H4").Value = 6.540496
$ws.Range("I4").Value = 0.01970539991828544
$ws.Range("J4").Value = 0.01970539991828544
$ws.Range("M4").Value = 5.413469333333334
$ws.Range("N4").Value = 16.240408
$ws.Range("O4").Value = 0.0525235026743817
$ws.Range("P4").Value = 0.0525235026743817
$ws.Range("Q4").Value = 11.80225817359645
$ws.Range("R4").Value = 106.220323562368
$ws.Range("S4").Value = 0.001034996625307826
$ws.Range("T4").Value = 0.001034996625307826
$ws.Range("G5").Value = 2.180165333333334
$ws.Range("H5").Value = 6.540496
$ws.Range("I5").Value = 0.01970539991828544
$ws.Range("J5").Value = 0.01970539991828544
$ws.Range("M5").Value = 47.79503400000001
$ws.Range("N5").Value = 143.385102
$ws.Range("O5").Value = 0.4637252825399149
$ws.Range("P5").Value = 0.4637252825399149
$ws.Range("Q5").Value = 104.201076232288
$ws.Range("R5").Value = 937.8096860905921
$ws.Range("S5").Value = 0.00913789214466893
$ws.Range("T5").Value = 0.00913789214466893
$ws.Range("I6").Value = 0.733713204346044
$ws.Range("J6").Value = 0.7337132043460441
$ws.Range("M6").Value = 28.31444233333334
$ws.Range("N6").Value = 84.94332700000001
$ws.Range("O6").Value = 0.2747173016130739
$ws.Range("P6").Value = 0.2747173016130739
$ws.Range("Q6").Value = 2298.468328515463
$ws.Range("R6").Value = 20686.21495663916
$ws.Range("S6").Value = 0.2015637116558271
$ws.Range("T6").Value = 0.2015637116558272
$ws.Range("I7").Value = 0.733713204346044
$ws.Range("J7").Value = 0.7337132043460441
$ws.Range("O7").Value = 0.2090339131726295
$ws.Range("P7").Value = 0.2090339131726295
$ws.Range("S7").Value = 0.1533709422508827
$ws.Range("T7").Value = 0.1533709422508828
$ws.Range("I8").Value = 0.733713204346044
$ws.Range("J8").Value = 0.7337132043460441
$ws.Range("M8").Value = 5.413469333333334
$ws.Range("N8").Value = 16.240408
$ws.Range("O8").Value = 0.0525235026743817
$ws.Range("P8").Value = 0.0525235026743817
$ws.Range("Q8").Value = 439.4466846132499
$ws.Range("R8").Value = 3955.020161519249
$ws.Range("S8").Value = 0.03853718745069861
$ws.Range("T8").Value = 0.03853718745069862
$ws.Range("I9").Value = 0.733713204346044
$ws.Range("J9").Value = 0.7337132043460441
$ws.Range("M9").Value = 47.79503400000001
$ws.Range("N9").Value = 143.385102
$ws.Range("O9").Value = 0.4637252825399149
$ws.Range("P9").Value = 0.4637252825399149
$ws.Range("Q9").Value = 3879.835266258869
$ws.Range("R9").Value = 34918.51739632981
$ws.Range("S9").Value = 0.3402413629886356
$ws.Range("T9").Value = 0.3402413629886356
$ws.Range("G10").Value = 25.672264
$ws.Range("H10").Value = 77.016792
$ws.Range("I10").Value = 0.2320384702908474
$ws.Range("J10").Value = 0.2320384702908474
$ws.Range("M10").Value = 28.31444233333334
$ws.Range("N10").Value = 84.94332700000001
$ws.Range("O10").Value = 0.2747173016130739
$ws.Range("P10").Value = 0.2747173016130739
$ws.Range("Q10").Value = 726.8958385941094
$ws.Range("R10").Value = 6542.062547346985
$ws.Range("S10").Value = 0.06374498242872702
$ws.Range("T10").Value = 0.06374498242872702
$ws.Range("G11").Value = 25.672264
$ws.Range("H11").Value = 77.016792
$ws.Range("I11").Value = 0.2320384702908474
$ws.Range("J11").Value = 0.2320384702908474
$ws.Range("O11").Value = 0.2090339131726295
$ws.Range("P11").Value = 0.2090339131726295
$ws.Range("Q11").Value = 553.0990611732027
$ws.Range("R11").Value = 4977.891550558824
$ws.Range("S11").Value = 0.04850390945148677
$ws.Range("T11").Value = 0.04850390945148677
$ws.Range("G12").Value = 25.672264
$ws.Range("H12").Value = 77.016792
$ws.Range("I12").Value = 0.2320384702908474
$ws.Range("J12").Value = 0.2320384702908474
$ws.Range("M12").Value = 5.413469333333334
$ws.Range("N12").Value = 16.240408
$ws.Range("O12").Value = 0.0525235026743817
$ws.Range("P12").Value = 0.0525235026743817
$ws.Range("Q12").Value = 138.9760138812373
$ws.Range("R12").Value = 1250.784124931136
$ws.Range("S12").Value = 0.01218747321488076
$ws.Range("T12").Value = 0.01218747321488076
$ws.Range("G13").Value = 25.672264
$ws.Range("H13").Value = 77.016792
$ws.Range("I13").Value = 0.2320384702908474
$ws.Range("J13").Value = 0.2320384702908474
$ws.Range("M13").Value = 47.79503400000001
$ws.Range("N13").Value = 143.385102
$ws.Range("O13").Value = 0.4637252825399149
$ws.Range("P13").Value = 0.4637252825399149
$ws.Range("Q13").Value = 1227.006730736976
$ws.Range("R13").Value = 11043.06057663278
$ws.Range("S13").Value = 0.1076021051957528
$ws.Range("T13").Value = 0.1076021051957528
$ws.Range("G14").Value = 1.608999666666667
$ws.Range("H14").Value = 4.826999
$ws.Range("I14").Value = 0.01454292544482312
$ws.Range("J14").Value = 0.01454292544482312
$ws.Range("M14").Value = 28.31444233333334
$ws.Range("N14").Value = 84.94332700000001
$ws.Range("O14").Value = 0.2747173016130739
$ws.Range("P14").Value = 0.2747173016130739
$ws.Range("Q14").Value = 45.55792827618589
$ws.Range("R14").Value = 410.021354485673
$ws.Range("S14").Value = 0.003995193235761922
$ws.Range("T14").Value = 0.003995193235761922
$ws.Range("G15").Value = 1.608999666666667
$ws.Range("H15").Value = 4.826999
$ws.Range("I15").Value = 0.01454292544482312
$ws.Range("J15").Value = 0.01454292544482312
$ws.Range("O15").Value = 0.2090339131726295
$ws.Range("P15").Value = 0.2090339131726295
$ws.Range("Q15").Value = 34.66527942612812
$ws.Range("R15").Value = 311.987514835153
$ws.Range("S15").Value = 0.003039964614709182
$ws.Range("T15").Value = 0.003039964614709182
$ws.Range("G16").Value = 1.608999666666667
$ws.Range("H16").Value = 4.826999
$ws.Range("I16").Value = 0.01454292544482312
$ws.Range("J16").Value = 0.01454292544482312
$ws.Range("M16").Value = 5.413469333333334
$ws.Range("N16").Value = 16.240408
$ws.Range("O16").Value = 0.0525235026743817
$ws.Range("P16").Value = 0.0525235026743817
$ws.Range("Q16").Value = 8.710270352843557
$ws.Range("R16").Value = 78.392433175592
$ws.Range("S16").Value = 0.0007638453834945011
$ws.Range("T16").Value = 0.0007638453834945011
$ws.Range("G17").Value = 1.608999666666667
$ws.Range("H17").Value = 4.826999
$ws.Range("I17").Value = 0.01454292544482312
$ws.Range("J17").Value = 0.01454292544482312
$ws.Range("M17").Value = 47.79503400000001
$ws.Range("N17").Value = 143.385102
$ws.Range("O17").Value = 0.4637252825399149
$ws.Range("P17").Value = 0.4637252825399149
$ws.Range("Q17").Value = 76.90219377432201
$ws.Range("R17").Value = 692.1197439688981
$ws.Range("S17").Value = 0.00674392221085752
$ws.Range("T17").Value = 0.00674392221085752

Write-Host "Applied 174 changes"